# cambios en PC Gestion Documental para guardar en R
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the phone/reference number in E2
$ws.Range("E2").Value = 9498924883

# Move the active selection to G3 (matches the saved cursor position)
$ws.Range("G3").Select()
